$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.20528676305393
$ws.Range("C2").Value = 10.30178542906352
$ws.Range("D2").Value = 4.433521130763651
$ws.Range("F2").Value = 23.26322633026942
$ws.Range("G2").Value = 3.617448312761777
$ws.Range("I2").Value = 20.97475683152234
$ws.Range("L2").Value = 10.65195905747185
$ws.Range("N2").Value = 16.89179742351937
$ws.Range("O2").Value = 20.68331587780488
$ws.Range("B3").Value = 14.65848222874885
$ws.Range("C3").Value = 10.10787278854541
$ws.Range("D3").Value = 4.395257537804144
$ws.Range("F3").Value = 23.21544168812196
$ws.Range("G3").Value = 3.619525761940682
$ws.Range("I3").Value = 21.06394726325459
$ws.Range("L3").Value = 10.62380305732493
$ws.Range("N3").Value = 16.93888867350491
$ws.Range("O3").Value = 20.70854798102424
$ws.Range("B4").Value = 14.31396454376042
$ws.Range("C4").Value = 9.985920972977109
$ws.Range("D4").Value = 4.371285496611651
$ws.Range("F4").Value = 23.19350196777786
$ws.Range("G4").Value = 3.620869760474023
$ws.Range("I4").Value = 21.12371969973346
$ws.Range("L4").Value = 10.60867511328675
$ws.Range("N4").Value = 16.9695925348456
$ws.Range("O4").Value = 20.72976496713104
$ws.Range("B5").Value = 14.17157772463309
$ws.Range("C5").Value = 9.935539730537435
$ws.Range("D5").Value = 4.361400684938728
$ws.Range("F5").Value = 23.18642824169083
$ws.Range("G5").Value = 3.621434714200854
$ws.Range("I5").Value = 21.14933359963494
$ws.Range("L5").Value = 10.60305800535437
$ws.Range("N5").Value = 16.98255544661242
$ws.Range("O5").Value = 20.73984705217921
$ws.Range("B6").Value = 14.14782079340942
$ws.Range("C6").Value = 9.927133845708983
$ws.Range("D6").Value = 4.359752436210472
$ws.Range("F6").Value = 23.18536654379675
$ws.Range("G6").Value = 3.62152956860939
$ws.Range("I6").Value = 21.15366253514701
$ws.Range("L6").Value = 10.60215848718624
$ws.Range("N6").Value = 16.98473518217748
$ws.Range("O6").Value = 20.74160779068901
$ws.Range("B7").Value = 14.31205203688453
$ws.Range("C7").Value = 9.985244230868958
$ws.Range("D7").Value = 4.371152650740159
$ws.Range("F7").Value = 23.19339900375135
$ws.Range("G7").Value = 3.620877309665768
$ws.Range("I7").Value = 21.12406005624813
$ws.Range("L7").Value = 10.60859713601166
$ws.Range("N7").Value = 16.96976553068649
$ws.Range("O7").Value = 20.7298951284794
$ws.Range("B8").Value = 15.01869233131083
$ws.Range("C8").Value = 10.2355482165308
$ws.Range("D8").Value = 4.420429561899932
$ws.Range("F8").Value = 23.24521854499897
$ws.Range("G8").Value = 3.618150446150101
$ws.Range("I8").Value = 21.004467982361
$ws.Range("L8").Value = 10.64180543509614
$ws.Range("N8").Value = 16.90766353550681
$ws.Range("O8").Value = 20.69082594342618
$ws.Range("B9").Value = 16.3264091609242
$ws.Range("C9").Value = 10.70174766796471
$ws.Range("D9").Value = 4.513087879004767
$ws.Range("F9").Value = 23.4051883277727
$ws.Range("G9").Value = 3.613343599648699
$ws.Range("I9").Value = 20.80984642454262
$ws.Range("L9").Value = 10.72384385163468
$ws.Range("N9").Value = 16.8000451340643
$ws.Range("O9").Value = 20.65975635892505
$ws.Range("B10").Value = 17.22985074187595
$ws.Range("C10").Value = 11.02707214853944
$ws.Range("D10").Value = 4.578505722302741
$ws.Range("F10").Value = 23.55765625118628
$ws.Range("G10").Value = 3.610138048994385
$ws.Range("I10").Value = 20.69138800291504
$ws.Range("L10").Value = 10.79410169628534
$ws.Range("N10").Value = 16.72956216320888
$ws.Range("O10").Value = 20.66482126284693
$ws.Range("B11").Value = 17.62670628818144
$ws.Range("C11").Value = 11.17092036406003
$ws.Range("D11").Value = 4.607640153047183
$ws.Range("F11").Value = 23.63442395687386
$ws.Range("G11").Value = 3.608749817637116
$ws.Range("I11").Value = 20.64286837888357
$ws.Range("L11").Value = 10.82815138209877
$ws.Range("N11").Value = 16.69935133808284
$ws.Range("O11").Value = 20.67319285881358
$ws.Range("B12").Value = 17.77483568159357
$ws.Range("C12").Value = 11.22476403366744
$ws.Range("D12").Value = 4.618578779303618
$ws.Range("F12").Value = 23.66454087892712
$ws.Range("G12").Value = 3.608234138862704
$ws.Range("I12").Value = 20.62527071285842
$ws.Range("L12").Value = 10.84133796441634
$ws.Range("N12").Value = 16.68817689810015
$ws.Range("O12").Value = 20.67723496282575
$ws.Range("B13").Value = 17.74303079115117
$ws.Range("C13").Value = 11.21319628744397
$ws.Range("D13").Value = 4.616227193974598
$ws.Range("F13").Value = 23.65800843068612
$ws.Range("G13").Value = 3.608344754917984
$ws.Range("I13").Value = 20.62902612751634
$ws.Range("L13").Value = 10.83848510703977
$ws.Range("N13").Value = 16.690571704048
$ws.Range("O13").Value = 20.67632566024677
$ws.Range("B14").Value = 17.63893672908714
$ws.Range("C14").Value = 11.17536289134043
$ws.Range("D14").Value = 4.608541986758891
$ws.Range("F14").Value = 23.63688083224741
$ws.Range("G14").Value = 3.608707192017851
$ws.Range("I14").Value = 20.64140504135605
$ws.Range("L14").Value = 10.8292304295568
$ws.Range("N14").Value = 16.69842668766271
$ws.Range("O14").Value = 20.67350793582356
$ws.Range("B15").Value = 17.57489265983812
$ws.Range("C15").Value = 11.15210605772602
$ws.Range("D15").Value = 4.603822218730145
$ws.Range("F15").Value = 23.6240752824391
$ws.Range("G15").Value = 3.608930497885084
$ws.Range("I15").Value = 20.64908860926821
$ws.Range("L15").Value = 10.82359955530862
$ws.Range("N15").Value = 16.70327268205981
$ws.Range("O15").Value = 20.67189552300573
$ws.Range("B16").Value = 17.20362045596151
$ws.Range("C16").Value = 11.01758510666169
$ws.Range("D16").Value = 4.576588795666877
$ws.Range("F16").Value = 23.5527868344234
$ws.Range("G16").Value = 3.610230177256275
$ws.Range("I16").Value = 20.69466721194658
$ws.Range("L16").Value = 10.79191785834601
$ws.Range("N16").Value = 16.73157373109869
$ws.Range("O16").Value = 20.66439622661946
$ws.Range("B17").Value = 16.97214892455167
$ws.Range("C17").Value = 10.93397565252696
$ws.Range("D17").Value = 4.559719046345382
$ws.Range("F17").Value = 23.51093864687859
$ws.Range("G17").Value = 3.611045378707546
$ws.Range("I17").Value = 20.72400572315453
$ws.Range("L17").Value = 10.77301188334006
$ws.Range("N17").Value = 16.74940947149239
$ws.Range("O17").Value = 20.66134957491216
$ws.Range("B18").Value = 16.83768929408067
$ws.Range("C18").Value = 10.8854975145716
$ws.Range("D18").Value = 4.549957476831582
$ws.Range("F18").Value = 23.48756732773759
$ws.Range("G18").Value = 3.611520851755873
$ws.Range("I18").Value = 20.74138546919397
$ws.Range("L18").Value = 10.76233478665404
$ws.Range("N18").Value = 16.75984249315107
$ws.Range("O18").Value = 20.66016841303761
$ws.Range("B19").Value = 16.7919402055002
$ws.Range("C19").Value = 10.86901800960049
$ws.Range("D19").Value = 4.546642451567086
$ws.Range("F19").Value = 23.47977473433139
$ws.Range("G19").Value = 3.611682972259236
$ws.Range("I19").Value = 20.74735657213689
$ws.Range("L19").Value = 10.75875378791044
$ws.Range("N19").Value = 16.76340490299743
$ws.Range("O19").Value = 20.65986659747325
$ws.Range("B20").Value = 16.99692734261043
$ws.Range("C20").Value = 10.94291644172074
$ws.Range("D20").Value = 4.561520948571099
$ws.Range("F20").Value = 23.51532127980909
$ws.Range("G20").Value = 3.610957917369672
$ws.Range("I20").Value = 20.72083029419735
$ws.Range("L20").Value = 10.77500410980616
$ws.Range("N20").Value = 16.74749278371101
$ws.Range("O20").Value = 20.66161478124641
$ws.Range("B21").Value = 17.66957094230913
$ws.Range("C21").Value = 11.18649278471553
$ws.Range("D21").Value = 4.610801900808721
$ws.Range("F21").Value = 23.64305827359942
$ws.Range("G21").Value = 3.608600464068305
$ws.Range("I21").Value = 20.63774797308947
$ws.Range("L21").Value = 10.8319408722466
$ws.Range("N21").Value = 16.69611228106403
$ws.Range("O21").Value = 20.67431191348195
$ws.Range("B22").Value = 18.09659844381068
$ws.Range("C22").Value = 11.3420068279014
$ws.Range("D22").Value = 4.64246038034022
$ws.Range("F22").Value = 23.73263241660718
$ws.Range("G22").Value = 3.607118080152247
$ws.Range("I22").Value = 20.58797174616848
$ws.Range("L22").Value = 10.87085457440242
$ws.Range("N22").Value = 16.66408084614622
$ws.Range("O22").Value = 20.68769183969404
$ws.Range("B23").Value = 17.86987304942739
$ws.Range("C23").Value = 11.25935275898681
$ws.Range("D23").Value = 4.625615277265003
$ws.Range("F23").Value = 23.68427443088498
$ws.Range("G23").Value = 3.607903933950115
$ws.Range("I23").Value = 20.61412317423703
$ws.Range("L23").Value = 10.84993249911189
$ws.Range("N23").Value = 16.68103511420959
$ws.Range("O23").Value = 20.68008615024548
$ws.Range("B24").Value = 16.98572932000375
$ws.Range("C24").Value = 10.93887558535538
$ws.Range("D24").Value = 4.560706504252862
$ws.Range("F24").Value = 23.51333774771131
$ws.Range("G24").Value = 3.610997437446597
$ws.Range("I24").Value = 20.72226430903105
$ws.Range("L24").Value = 10.77410282496451
$ws.Range("N24").Value = 16.7483587605701
$ws.Range("O24").Value = 20.66149310465155
$ws.Range("B25").Value = 15.98207282687133
$ws.Range("C25").Value = 10.57849868768467
$ws.Range("D25").Value = 4.488471344762172
$ws.Range("F25").Value = 23.35572217997017
$ws.Range("G25").Value = 3.61458647319513
$ws.Range("I25").Value = 20.85820381731295
$ws.Range("L25").Value = 10.69987112498742
$ws.Range("N25").Value = 16.8276477385981
$ws.Range("O25").Value = 20.66326950348931
